$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add pre-treatment phase parent outcome measure values in column C
$ws.Range("C2").Value = "A little stressful"
$ws.Range("C3").Value = "A little stressful"
$ws.Range("C4").Value = "Moderately stressful"
$ws.Range("C5").Value = "Not stressful"
$ws.Range("C6").Value = "Moderately stressful"
$ws.Range("C7").Value = "Moderately stressful"

# Update the active cell selection
$ws.Range("C8").Select()
